# Update cryptocurrency price/volume data in the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "33.990.71"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.782.42"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.23"
$ws.Range("E5").Value = "  +2.52%  "
$ws.Range("E6").Value = "  +1.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.33"
$ws.Range("E8").Value = "  +3.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.292"
$ws.Range("E9").Value = "  +2.65%  "
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0939"
$ws.Range("E11").Value = "  +1.90%  "
$ws.Range("D12").Value = "2.037.72"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.10"
$ws.Range("E13").Value = "  +5.67%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.789.01"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D16").Value = "33.982.77"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.17"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.89"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.21"
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").Value = "0.0{0}0785" -f [char]0x2083
$ws.Range("E20").Value = "  +2.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.73"
$ws.Range("E22").Value = "  +2.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.12"
$ws.Range("E23").Value = "  +2.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.06"
$ws.Range("E24").Value = "  -2.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.01"
$ws.Range("E25").Value = "  +1.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.32"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  +1.91%  "
$ws.Range("E28").Value = "  +1.65%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  +3.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0513"
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("E33").Value = "  +1.29%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "1.396.60"
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.665"
$ws.Range("E36").Value = "  +5.89%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("E38").Value = "  +1.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.25"
$ws.Range("E39").Value = "  +7.33%  "
$ws.Range("E40").Value = "  +0.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.915"
$ws.Range("E41").Value = "  -1.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "78.11"
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.41"
$ws.Range("E44").Value = "  +13.80%  "
$ws.Range("D45").Value = "0.0{0}0145" -f [char]0x2086
$ws.Range("E45").Value = "  +23.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.08"
$ws.Range("E46").Value = "  +4.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "108.16"
$ws.Range("E47").Value = "  +4.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0497"
$ws.Range("E48").Value = "  +1.71%  "
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").Value = "1.937.17"
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("E51").Value = "  +0.47%  "
